$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 2
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 0

$ws.Range("A1:B5").Font.Name = "Arial"

$ws.Range("C7").Select()
